$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 74.609651
$ws.Range("H2").Value = 223.828953
$ws.Range("I2").Value = 0.1061386348809139
$ws.Range("J2").Value = 0.1061386348809139
$ws.Range("M2").Value = 2.598166333333333
$ws.Range("N2").Value = 7.794499
$ws.Range("O2").Value = 0.3466013321552429
$ws.Range("P2").Value = 0.3466013321552429
$ws.Range("Q2").Value = 193.8482833699497
$ws.Range("R2").Value = 1744.634550329547
$ws.Range("S2").Value = 0.03678779224286369
$ws.Range("T2").Value = 0.0367877922428637
$ws.Range("G3").Value = 74.609651
$ws.Range("H3").Value = 223.828953
$ws.Range("I3").Value = 0.1061386348809139
$ws.Range("J3").Value = 0.1061386348809139
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("O3").Value = 0.5780859172985858
$ws.Range("P3").Value = 0.5780859172985858
$ws.Range("Q3").Value = 323.3137103422366
$ws.Range("R3").Value = 2909.82339308013
$ws.Range("S3").Value = 0.0613572501059528
$ws.Range("T3").Value = 0.06135725010595281
$ws.Range("G4").Value = 74.609651
$ws.Range("H4").Value = 223.828953
$ws.Range("I4").Value = 0.1061386348809139
$ws.Range("J4").Value = 0.1061386348809139
$ws.Range("M4").Value = 0.4692043333333333
$ws.Range("N4").Value = 1.407613
$ws.Range("O4").Value = 0.06259293136852516
$ws.Range("P4").Value = 0.06259293136852516
$ws.Range("Q4").Value = 35.00717155768767
$ws.Range("R4").Value = 315.064544019189
$ws.Range("S4").Value = 0.006643528288649995
$ws.Range("T4").Value = 0.006643528288649996
$ws.Range("G5").Value = 74.609651
$ws.Range("H5").Value = 223.828953
$ws.Range("I5").Value = 0.1061386348809139
$ws.Range("J5").Value = 0.1061386348809139
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09534933333333333
$ws.Range("N5").Value = 0.286048
$ws.Range("O5").Value = 0.01271981917764605
$ws.Range("P5").Value = 0.01271981917764604
$ws.Range("Q5").Value = 7.113980483082666
$ws.Range("R5").Value = 64.02582434774399
$ws.Range("S5").Value = 0.00135006424344742
$ws.Range("T5").Value = 0.00135006424344742
$ws.Range("G6").Value = 597.374756
$ws.Range("I6").Value = 0.8498168837991085
$ws.Range("J6").Value = 0.8498168837991086
$ws.Range("M6").Value = 2.598166333333333
$ws.Range("N6").Value = 7.794499
$ws.Range("O6").Value = 0.3466013321552429
$ws.Range("P6").Value = 0.3466013321552429
$ws.Range("Q6").Value = 1552.078979422415
$ws.Range("R6").Value = 13968.71081480173
$ws.Range("S6").Value = 0.2945476640127883
$ws.Range("T6").Value = 0.2945476640127883
$ws.Range("G7").Value = 597.374756
$ws.Range("I7").Value = 0.8498168837991085
$ws.Range("J7").Value = 0.8498168837991086
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("O7").Value = 0.5780859172985858
$ws.Range("P7").Value = 0.5780859172985858
$ws.Range("S7").Value = 0.4912671728068334
$ws.Range("T7").Value = 0.4912671728068334
$ws.Range("G8").Value = 597.374756
$ws.Range("I8").Value = 0.8498168837991085
$ws.Range("J8").Value = 0.8498168837991086
$ws.Range("M8").Value = 0.4692043333333333
$ws.Range("N8").Value = 1.407613
$ws.Range("O8").Value = 0.06259293136852516
$ws.Range("P8").Value = 0.06259293136852516
$ws.Range("Q8").Value = 280.2908241391427
$ws.Range("R8").Value = 2522.617417252284
$ws.Range("S8").Value = 0.05319252988345152
$ws.Range("T8").Value = 0.05319252988345152
$ws.Range("G9").Value = 597.374756
$ws.Range("I9").Value = 0.8498168837991085
$ws.Range("J9").Value = 0.8498168837991086
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09534933333333333
$ws.Range("N9").Value = 0.286048
$ws.Range("O9").Value = 0.01271981917764605
$ws.Range("P9").Value = 0.01271981917764604
$ws.Range("Q9").Value = 56.95928473476267
$ws.Range("R9").Value = 512.6335626128639
$ws.Range("S9").Value = 0.0108095170960353
$ws.Range("T9").Value = 0.0108095170960353
$ws.Range("G10").Value = 30.48438
$ws.Range("H10").Value = 91.45313999999999
$ws.Range("I10").Value = 0.04336664808137267
$ws.Range("J10").Value = 0.04336664808137267
$ws.Range("M10").Value = 2.598166333333333
$ws.Range("N10").Value = 7.794499
$ws.Range("O10").Value = 0.3466013321552429
$ws.Range("P10").Value = 0.3466013321552429
$ws.Range("Q10").Value = 79.20348980854
$ws.Range("R10").Value = 712.83140827686
$ws.Range("S10").Value = 0.01503093799611138
$ws.Range("T10").Value = 0.01503093799611138
$ws.Range("G11").Value = 30.48438
$ws.Range("H11").Value = 91.45313999999999
$ws.Range("I11").Value = 0.04336664808137267
$ws.Range("J11").Value = 0.04336664808137267
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("O11").Value = 0.5780859172985858
$ws.Range("P11").Value = 0.5780859172985858
$ws.Range("Q11").Value = 132.1011139066
$ws.Range("R11").Value = 1188.9100251594
$ws.Range("S11").Value = 0.02506964853628528
$ws.Range("T11").Value = 0.02506964853628528
$ws.Range("G12").Value = 30.48438
$ws.Range("H12").Value = 91.45313999999999
$ws.Range("I12").Value = 0.04336664808137267
$ws.Range("J12").Value = 0.04336664808137267
$ws.Range("M12").Value = 0.4692043333333333
$ws.Range("N12").Value = 1.407613
$ws.Range("O12").Value = 0.06259293136852516
$ws.Range("P12").Value = 0.06259293136852516
$ws.Range("Q12").Value = 14.30340319498
$ws.Range("R12").Value = 128.73062875482
$ws.Range("S12").Value = 0.002714445627040343
$ws.Range("T12").Value = 0.002714445627040343
$ws.Range("G13").Value = 30.48438
$ws.Range("H13").Value = 91.45313999999999
$ws.Range("I13").Value = 0.04336664808137267
$ws.Range("J13").Value = 0.04336664808137267
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09534933333333333
$ws.Range("N13").Value = 0.286048
$ws.Range("O13").Value = 0.01271981917764605
$ws.Range("P13").Value = 0.01271981917764604
$ws.Range("Q13").Value = 2.90666531008
$ws.Range("R13").Value = 26.15998779072
$ws.Range("S13").Value = 0.0005516159219356711
$ws.Range("T13").Value = 0.000551615921935671
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.4764796666666666
$ws.Range("H14").Value = 1.429439
$ws.Range("I14").Value = 0.0006778332386049212
$ws.Range("J14").Value = 0.0006778332386049213
$ws.Range("M14").Value = 2.598166333333333
$ws.Range("N14").Value = 7.794499
$ws.Range("O14").Value = 0.3466013321552429
$ws.Range("P14").Value = 0.3466013321552429
$ws.Range("Q14").Value = 1.237973428451222
$ws.Range("R14").Value = 11.141760856061
$ws.Range("S14").Value = 0.0002349379034795683
$ws.Range("T14").Value = 0.0002349379034795683
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.4764796666666666
$ws.Range("H15").Value = 1.429439
$ws.Range("I15").Value = 0.0006778332386049212
$ws.Range("J15").Value = 0.0006778332386049213
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("O15").Value = 0.5780859172985858
$ws.Range("P15").Value = 0.5780859172985858
$ws.Range("Q15").Value = 2.064778575798889
$ws.Range("R15").Value = 18.58300718219
$ws.Range("S15").Value = 0.0003918458495143971
$ws.Range("T15").Value = 0.0003918458495143971
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.4764796666666666
$ws.Range("H16").Value = 1.429439
$ws.Range("I16").Value = 0.0006778332386049212
$ws.Range("J16").Value = 0.0006778332386049213
$ws.Range("M16").Value = 0.4692043333333333
$ws.Range("N16").Value = 1.407613
$ws.Range("O16").Value = 0.06259293136852516
$ws.Range("P16").Value = 0.06259293136852516
$ws.Range("Q16").Value = 0.2235663243452222
$ws.Range("R16").Value = 2.012096919107
$ws.Range("S16").Value = 0.00004242756938330297
$ws.Range("T16").Value = 0.00004242756938330298
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.4764796666666666
$ws.Range("H17").Value = 1.429439
$ws.Range("I17").Value = 0.0006778332386049212
$ws.Range("J17").Value = 0.0006778332386049213
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09534933333333333
$ws.Range("N17").Value = 0.286048
$ws.Range("O17").Value = 0.01271981917764605
$ws.Range("P17").Value = 0.01271981917764604
$ws.Range("Q17").Value = 0.04543201856355555
$ws.Range("R17").Value = 0.4088881670719999
$ws.Range("S17").Value = 0.000008621916227652805
$ws.Range("T17").Value = 0.000008621916227652805
